# This script updates the "北京-漫展信息" workbook to reflect the latest
# scrape (commit: "Update gh-pages to output generated at 456a3b4").
#
# Changes:
#  1) Numeric "想去人数" (want-to-go count, column F) updates across the
#     "展览", "演出", "本地生活" and "全部类型" sheets.
#  2) On the "全部类型" sheet, the event that used to occupy row 20
#     ("北京·第16届IJOY漫展【《大主宰年番》...专场见面会】") dropped out of the
#     combined listing, so rows 20-24 each take on the data previously held
#     by the row below them, and a new event ("北京·摇滚新星企划") is placed
#     into row 25.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) F-column (想去人数) updates ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 83
$ws1.Cells.Item(5, 6).Value = 9640
$ws1.Cells.Item(6, 6).Value = 643
$ws1.Cells.Item(8, 6).Value = 169
$ws1.Cells.Item(9, 6).Value = 317
$ws1.Cells.Item(10, 6).Value = 422
$ws1.Cells.Item(12, 6).Value = 200
$ws1.Cells.Item(13, 6).Value = 21
$ws1.Cells.Item(14, 6).Value = 463
$ws1.Cells.Item(15, 6).Value = 12279
$ws1.Cells.Item(20, 6).Value = 35
$ws1.Cells.Item(22, 6).Value = 46
$ws1.Cells.Item(26, 6).Value = 178
$ws1.Cells.Item(27, 6).Value = 160
$ws1.Cells.Item(28, 6).Value = 2733
$ws1.Cells.Item(31, 6).Value = 2104
$ws1.Cells.Item(34, 6).Value = 2153
$ws1.Cells.Item(35, 6).Value = 1030
$ws1.Cells.Item(36, 6).Value = 4214
$ws1.Cells.Item(37, 6).Value = 3682
$ws1.Cells.Item(38, 6).Value = 621
$ws1.Cells.Item(41, 6).Value = 775
$ws1.Cells.Item(42, 6).Value = 30
$ws1.Cells.Item(43, 6).Value = 118
$ws1.Cells.Item(44, 6).Value = 556
$ws1.Cells.Item(47, 6).Value = 239

# --- Sheet "演出" (performances) F-column (想去人数) updates ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(15, 6).Value = 25
$ws2.Cells.Item(17, 6).Value = 33
$ws2.Cells.Item(24, 6).Value = 79

# --- Sheet "本地生活" (local life) F-column (想去人数) updates ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 52

# --- Sheet "全部类型" (all types) F-column (想去人数) updates (non-shifted rows) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(6, 6).Value = 9640
$ws4.Cells.Item(7, 6).Value = 643
$ws4.Cells.Item(10, 6).Value = 169
$ws4.Cells.Item(11, 6).Value = 317
$ws4.Cells.Item(12, 6).Value = 422
$ws4.Cells.Item(13, 6).Value = 200
$ws4.Cells.Item(14, 6).Value = 21
$ws4.Cells.Item(15, 6).Value = 463
$ws4.Cells.Item(16, 6).Value = 12279
$ws4.Cells.Item(26, 6).Value = 178
$ws4.Cells.Item(27, 6).Value = 160
$ws4.Cells.Item(28, 6).Value = 2733
$ws4.Cells.Item(29, 6).Value = 2104
$ws4.Cells.Item(31, 6).Value = 2153
$ws4.Cells.Item(32, 6).Value = 1030
$ws4.Cells.Item(36, 6).Value = 4214
$ws4.Cells.Item(37, 6).Value = 3682
$ws4.Cells.Item(38, 6).Value = 621
$ws4.Cells.Item(41, 6).Value = 775
$ws4.Cells.Item(42, 6).Value = 30
$ws4.Cells.Item(43, 6).Value = 118
$ws4.Cells.Item(44, 6).Value = 556
$ws4.Cells.Item(47, 6).Value = 239

# --- Sheet "全部类型": rows 20-25 shifted up one (the row-20 "IJOY...孙郎朗"
#     event dropped off the combined feed) and a new row 25 ("摇滚新星企划") ---
$ws4.Cells.Item(20, 2).Value = "'2024-05-03"
$ws4.Cells.Item(20, 3).Value = "北京·塔罗集市"
$ws4.Cells.Item(20, 4).Value = "北京朝阳区广渠路南侧汇泰大厦1层展厅 汇泰大厦"
$ws4.Cells.Item(20, 5).Value = "2024.05.03 09:30-05.03 16:30"
$ws4.Cells.Item(20, 6).Value = 52
$ws4.Cells.Item(20, 7).Value = 75
$ws4.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84239"
$ws4.Cells.Item(20, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/z8qnfmoq1712735872200.jpeg"

$ws4.Cells.Item(21, 2).Value = "'2024-05-03"
$ws4.Cells.Item(21, 3).Value = "北京·知名演员 川久保拓司 专场活动"
$ws4.Cells.Item(21, 4).Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws4.Cells.Item(21, 5).Value = "2024.05.03 10:30-05.03 15:00"
$ws4.Cells.Item(21, 6).Value = 156
$ws4.Cells.Item(21, 7).Value = 528
$ws4.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82897"
$ws4.Cells.Item(21, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/rxrJuuvX1710409029498.jpeg"

$ws4.Cells.Item(22, 2).Value = "'2024-05-04"
$ws4.Cells.Item(22, 3).Value = "北京·XW咒术回战only"
$ws4.Cells.Item(22, 4).Value = "北花园路1号 超级蜂巢"
$ws4.Cells.Item(22, 5).Value = "2024.05.04 10:00-05.04 17:00"
$ws4.Cells.Item(22, 6).Value = 242
$ws4.Cells.Item(22, 7).Value = 60
$ws4.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83570"
$ws4.Cells.Item(22, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/G9X2HmU11711703284044.jpeg"

$ws4.Cells.Item(23, 2).Value = "'2024-05-04"
$ws4.Cells.Item(23, 3).Value = "北京·第16届IJOY漫展【文森个人专场见面会】"
$ws4.Cells.Item(23, 4).Value = "天辰东路7号 北京国家会议中心"
$ws4.Cells.Item(23, 5).Value = "2024.05.04 11:00-05.04 15:10"
$ws4.Cells.Item(23, 6).Value = 41
$ws4.Cells.Item(23, 7).Value = 238
$ws4.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83617"
$ws4.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/E3ZY4mKN1711961443069.jpeg"

$ws4.Cells.Item(24, 2).Value = "'2024-05-12"
$ws4.Cells.Item(24, 3).Value = "北京·《家庭教师》《七龙珠》《火影忍者》超燃动漫音乐会"
$ws4.Cells.Item(24, 4).Value = "北京东图剧场 北京东图剧场"
$ws4.Cells.Item(24, 5).Value = "2024.05.12 19:30-05.12 21:00"
$ws4.Cells.Item(24, 6).Value = 51
$ws4.Cells.Item(24, 7).Value = 64
$ws4.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84067"
$ws4.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/9gLpckTZ1712754110725.jpeg"

$ws4.Cells.Item(25, 2).Value = "'2024-05-12"
$ws4.Cells.Item(25, 3).Value = "北京·摇滚新星企划"
$ws4.Cells.Item(25, 4).Value = "朝阳北路甲27号菁英梦谷·常营文创产业园南门B5座 WeShow Live 北京"
$ws4.Cells.Item(25, 5).Value = "2024.05.12 15:00-05.12 17:50"
$ws4.Cells.Item(25, 6).Value = 25
$ws4.Cells.Item(25, 7).Value = 99
$ws4.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84069"
$ws4.Cells.Item(25, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/1Hz7WZo91712620004229.jpeg"

Write-Host "Applied want-to-go-count updates and refreshed combined-sheet rows 20-25."
